$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.829.10'
$ws.Range('E2').Value = '  -1.34%  '

# Row 3
$ws.Range('D3').Value = '1.873.34'
$ws.Range('E3').Value = '  -1.56%  '

# Row 4
$ws.Range('E4').Value = '  -0.09%  '

# Row 5
$ws.Range('D5').Value = "'301.68"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.89%  '

# Row 6
$ws.Range('E6').Value = '  -0.07%  '

# Row 7
$ws.Range('D7').Value = "'0.5379"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.49%  '

# Row 8
$ws.Range('D8').Value = "'0.3761"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.41%  '

# Row 9
$ws.Range('D9').Value = "'0.07185"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.60%  '

# Row 10
$ws.Range('D10').Value = "'21.56"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.22%  '

# Row 11
$ws.Range('D11').Value = "'0.8881"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.88%  '

# Row 12
$ws.Range('D12').Value = "'0.08155"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.28%  '

# Row 13
$ws.Range('D13').Value = '1.874.18'
$ws.Range('E13').Value = '  +0.47%  '

# Row 14
$ws.Range('D14').Value = "'93.43"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.27%  '

# Row 15
$ws.Range('D15').Value = "'5.261"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.72%  '

# Row 16
$ws.Range('D16').Value = "'1.002"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.07%  '

# Row 17
$ws.Range('D17').Value = "'14.74"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.18%  '

# Row 18
$ws.Range('D18').Value = "'0.000008550"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.29%  '

# Row 19
$ws.Range('E19').Value = '  -0.02%  '

# Row 20
$ws.Range('D20').Value = '26.890.42'
$ws.Range('E20').Value = '  -1.24%  '

# Row 21
$ws.Range('D21').Value = "'4.982"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.50%  '

# Row 22
$ws.Range('D22').Value = "'10.69"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.03%  '

# Row 23
$ws.Range('D23').Value = "'6.389"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.11%  '

# Row 24
$ws.Range('D24').Value = "'147.12"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.46%  '

# Row 25
$ws.Range('D25').Value = "'2.259"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.97%  '

# Row 26
$ws.Range('D26').Value = "'1.737"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.34%  '

# Row 27
$ws.Range('D27').Value = "'18.03"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.18%  '

# Row 28
$ws.Range('D28').Value = "'114.08"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.76%  '

# Row 29
$ws.Range('D29').Value = "'4.724"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.24%  '

# Row 30
$ws.Range('D30').Value = "'4.601"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.89%  '

# Row 31
$ws.Range('D31').Value = "'0.09159"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.85%  '

# Row 32
$ws.Range('D32').Value = "'0.8042"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.09%  '

# Row 33
$ws.Range('D33').Value = "'0.04975"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.94%  '

# Row 34
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = "'1.172"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.53%  '

# Row 35
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = "'2.981"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.11%  '

# Row 36
$ws.Range('D36').Value = "'0.6024"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.32%  '

# Row 37
$ws.Range('D37').Value = "'3.194"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.15%  '

# Row 38
$ws.Range('D38').Value = "'2.584"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.94%  '

# Row 39
$ws.Range('D39').Value = "'0.01954"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.86%  '

# Row 40
$ws.Range('D40').Value = "'1.073"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.98%  '

# Row 41
$ws.Range('D41').Value = "'6.573"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.25%  '

# Row 42
$ws.Range('D42').Value = "'8.842"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.10%  '

# Row 43
$ws.Range('B43').Value = 'Decentraland'
$ws.Range('C43').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D43').Value = "'0.5157"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.63%  '

# Row 44
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = "'115.43"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.70%  '

# Row 45
$ws.Range('D45').Value = "'0.1492"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.44%  '

# Row 46
$ws.Range('D46').Value = "'1.001"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.15%  '

# Row 47
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'9.963"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.19%  '

# Row 48
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = "'1.632"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.06%  '

# Row 49
$ws.Range('D49').Value = "'37.60"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.46%  '

# Row 50
$ws.Range('D50').Value = "'0.06027"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.15%  '

# Row 51
$ws.Range('D51').Value = "'62.15"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.91%  '
